# Draft profiles from the defunct "develop" branch
#
# Updates the "Data" worksheet (OutcomeOfCare - STU3 mapping notes) to
# reflect the latest draft mapping notes taken from the old "develop"
# branch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# P3: update the suggested FHIR resource mapping
$ws.Range("P3").Value = "CarePlan  / DiagnosticResult"

# R3: trim the long explanatory note down to just the heading
$ws.Range("R3").Value = "** OutcomeOfCare`n"

# P4: replace the DiagnosticReport mapping suggestion with the CarePlan one
$ws.Range("P4").Value = "CarePlan.activity:nursingIntervention.outcomeCodeableConcept Or derived profile on zib-TextResult."

# Q4: clear the stray "equal" equivalence note
$ws.Range("Q4").Value = ""

# P7: update the suggested FHIR resource mapping
$ws.Range("P7").Value = "Careplan.activity / DiagnosticReport.extention.partOf"

# R7: add a cautionary note about the extension choice
$ws.Range("R7").Value = "Maybe not the most suitable extension."

# Leave the cursor where the author last left it while drafting these notes
$ws.Activate()
$ws.Range("R13").Select() | Out-Null
